# Scene 54A edit:
#  1. Duplicate the full paragraph/table style block (Normal, TableNormal,
#     Heading1-6, Title, Subtitle) and insert the duplicate right before the
#     document's existing Subtitle style definition in word/styles.xml.
#  2. Add two new Google-Docs-style custom XML parts to the package:
#     customXML/item1.xml and customXML/itemProps1.xml.
#
# Both changes are made by round-tripping the package through
# Document.WordOpenXML (the Flat-OPC representation of the whole .docx),
# since they touch package-level structure (new parts) as well as the
# styles part.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# ---------------------------------------------------------------------
# 1. Duplicate the style block.
# ---------------------------------------------------------------------
$styleBlockStart = '<w:style w:type="paragraph" w:default="1" w:styleId="Normal">'
$subtitleMarker  = '<w:style w:type="paragraph" w:styleId="Subtitle">'
$stylesEnd       = '</w:styles>'

$blockStartIdx = $xml.IndexOf($styleBlockStart)
$stylesEndIdx  = $xml.IndexOf($stylesEnd)

# Full existing style set: Normal .. TableNormal .. Heading1-6 .. Title ..
# Subtitle (one full copy, as currently present in the document).
$allStylesBlock = $xml.Substring($blockStartIdx, $stylesEndIdx - $blockStartIdx)

# Insert a duplicate copy of that whole block immediately before the
# existing (single) Subtitle style definition, so the style set now
# appears twice in a row and the original Subtitle stays right after it.
$subtitleIdx = $xml.IndexOf($subtitleMarker)
$xml = $xml.Substring(0, $subtitleIdx) + $allStylesBlock + $xml.Substring($subtitleIdx)

# ---------------------------------------------------------------------
# 2. Add the two new customXML parts.
# ---------------------------------------------------------------------
$item1Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<go:gDocsCustomXmlDataStorage xmlns:go="http://customooxmlschemas.google.com/" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <go:docsCustomData xmlns:go="http://customooxmlschemas.google.com/" roundtripDataSignature="AMtx7mjBkRiTEwXbgkHdPIZhlrQcIpgbvg==">AMUW2mUla6W11047F3tF0Q2cVCxKw4Vdu8DljJaQx1J05P/93IXSuJXwZR6gbC7DJdmTmWVv7HEtmXgSoGunjgN2f71mn4IdDLVBlISz0lXPxyFDfCdiWu4=</go:docsCustomData>
</go:gDocsCustomXmlDataStorage>
'@

$itemProps1Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<ds:datastoreItem xmlns:ds="http://schemas.openxmlformats.org/officeDocument/2006/customXml" ds:itemID="{11111111-1234-1234-1234-123412341234}">
  <ds:schemaRefs>
    <ds:schemaRef ds:uri="http://schemas.openxmlformats.org/officeDocument/2006/relationships"/>
    <ds:schemaRef ds:uri="http://customooxmlschemas.google.com/"/>
  </ds:schemaRefs>
</ds:datastoreItem>
'@

$item1Part = '<pkg:part pkg:name="/customXml/item1.xml" pkg:contentType="application/xml"><pkg:xmlData>' + $item1Xml + '</pkg:xmlData></pkg:part>'
$itemProps1Part = '<pkg:part pkg:name="/customXml/itemProps1.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.customXmlProperties+xml"><pkg:xmlData>' + $itemProps1Xml + '</pkg:xmlData></pkg:part>'

$packageEnd = '</pkg:package>'
$packageEndIdx = $xml.IndexOf($packageEnd)
$xml = $xml.Substring(0, $packageEndIdx) + $item1Part + $itemProps1Part + $xml.Substring($packageEndIdx)

# ---------------------------------------------------------------------
# Write the rebuilt package back.
# ---------------------------------------------------------------------
$d.WordOpenXML = $xml
